$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "E-filing exemption - Circuit Court" (row 6) to "E-filing exemption".
#    Its URL (B6) is unchanged.
$ws.Range("A6").Value = "E-filing exemption"

# 2. Insert a new row for "Living will" at position 13 so the list stays
#    alphabetically sorted (pushes the old rows 13-15 down to 14-16).
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "Living will"
$ws.Range("B13").Value = "https://www.illinoislegalaid.org/legal-information/living-will"
$ws.Range("B13").Style = "Hyperlink"

# 3. The row insert above does not relocate the existing hyperlink objects,
#    so rebuild the hyperlinks collection from scratch to match the new
#    row layout (the new "Living will" row has no hyperlink, matching the
#    source data).
$ws.Range("B1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.illinoislegalaid.org/legal-information/appearance")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.illinoislegalaid.org/legal-information/fee-waiver")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand")
$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-circuit-court")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-appellate-court")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-supreme-court")

# Re-apply the Hyperlink cell style (Hyperlinks.Add can otherwise swap in a
# duplicate style record) so the styled cells match the original formatting.
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"
$ws.Range("B14").Style = "Hyperlink"
$ws.Range("B15").Style = "Hyperlink"
$ws.Range("B16").Style = "Hyperlink"
